# Update "想去人数" (interest count) values as part of the gh-pages data refresh.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 307
$wsExhibition.Range("F4").Value = 1242
$wsExhibition.Range("F5").Value = 621

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 307
$wsAll.Range("F4").Value = 1242
$wsAll.Range("F6").Value = 621
